$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132. This shifts the existing rows 132..177 down to
# 133..178 (carrying all their data and formatting with them), matching the
# diff's row-shift pattern. Excel's native Insert() also copies the row-above
# formatting into the freshly inserted row (e.g. the date style on column D).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with this week's data.
$ws.Range("A132").Value = 7
$ws.Range("B132").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C132").Value = "Ñuble"
$ws.Range("D132").Value = 44468
$ws.Range("E132").Value = 16
$ws.Range("F132").Value = 100114013
$ws.Range("G132").Value = "Zanahoria"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 160
$ws.Range("K132").Value = 6500
$ws.Range("L132").Value = 7000
$ws.Range("M132").Value = 6750
$ws.Range("N132").Value = "`$/saco 20 kilos"
$ws.Range("O132").Value = "Provincia de Diguillín"
$ws.Range("P132").Value = 338
$ws.Range("Q132").Value = 20
$ws.Range("R132").Value = "Hortaliza"
